$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) retains text formatting so values like "1.60" or "0.999"
# are not auto-converted to numbers by Excel when assigned as strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.752.52'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '2.632.86'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '579.78'
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("D6").Value = '155.17'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.622'
$ws.Range("E8").Value = '  -3.92%  '
$ws.Range("D9").Value = '2.630.34'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  -3.58%  '
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("E12").Value = '  -1.81%  '
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("D14").Value = '28.42'
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").Value = '3.105.56'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("D17").Value = '63.669.56'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '2.654.44'
$ws.Range("E18").Value = '  +2.37%  '
$ws.Range("D19").Value = '12.13'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").Value = '7.64'
$ws.Range("E20").Value = '  +3.80%  '
$ws.Range("D21").Value = '4.52'
$ws.Range("E21").Value = '  -3.30%  '
$ws.Range("D22").Value = '344.58'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '1.92'
$ws.Range("E24").Value = '  +10.48%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '67.88'
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("D26").Value = '0.0000109'
$ws.Range("E26").Value = '  -3.16%  '
$ws.Range("D27").Value = '599.05'
$ws.Range("E27").Value = '  +8.26%  '
$ws.Range("D28").Value = '9.26'
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").Value = '1.60'
$ws.Range("E29").Value = '  +3.23%  '
$ws.Range("D30").Value = '8.08'
$ws.Range("E30").Value = '  +2.30%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '1.75'
$ws.Range("E34").Value = '  +1.74%  '
$ws.Range("D35").Value = '6.61'
$ws.Range("E35").Value = '  +3.32%  '
$ws.Range("D36").Value = '5.44'
$ws.Range("E36").Value = '  +2.94%  '
$ws.Range("D37").Value = '0.403'
$ws.Range("E37").Value = '  -2.21%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").Value = '19.73'
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("D41").Value = '149.92'
$ws.Range("E41").Value = '  -2.16%  '
$ws.Range("E43").Value = '  +4.21%  '
$ws.Range("D44").Value = '41.87'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '159.34'
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '24.65'
$ws.Range("E46").Value = '  +8.47%  '
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").Value = '0.0588'
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("D49").Value = '0.632'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = '0.0999'
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("E51").Value = '  -0.52%  '
